$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) and Average (column D) with new Madigan bike hours data
$ws.Range("C2").Value = 99
$ws.Range("D2").Value = 226.75

$ws.Range("C3").Value = 185
$ws.Range("D3").Value = 211.64

$ws.Range("C4").Value = 221
$ws.Range("D4").Value = 212.32

$ws.Range("C5").Value = 264
$ws.Range("D5").Value = 239.33

$ws.Range("C6").Value = 260
$ws.Range("D6").Value = 241.59

$ws.Range("C7").Value = 78
$ws.Range("D7").Value = 113.83

$ws.Range("C8").Value = 74
$ws.Range("D8").Value = 94.23999999999999
